$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.240.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.337.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "411.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.330.83"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.629"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.115"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "39.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.72%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.868.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.333.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.056.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.35%  "
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000118"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "293.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -1.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.36%  "
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.83%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +18.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0490"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "134.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.155.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.670.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
